$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update mutex bit values (columns E-J hold bits; column C is a shared
# formula that sums them with weights 32,16,8,4,2,1)
$ws.Range("G4").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("G12").Value = 0

# Recalculate so the shared formula results in column C refresh
$excel.Calculate()

# Update the selected cell to reflect the last edited cell
$ws.Activate()
$ws.Range("G12").Select()
